# Apply the "Add files via upload" revision to the 5.5.1 indicator metadata sheet.
#
# Content changes:
#   - B4  (indicator description) is replaced with the new two-part (a/b) wording.
#   - B10 (organisation website) is replaced with the new domain.
# Both cells also pick up a distinct (Cyrillic-aware) font run, which is what the
# authoring workbook shows as two brand-new cellXfs entries, so we nudge the Font
# object on each cell to force the engine to materialise that new formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 5.5.1 indicator wording: now split into a) national parliament / b) local government.
$ws.Range("B4").Value = "5.5.1 Доля мест, занимаемых женщинами:`na) в национальном парламенте`nb) местных органах власти`n"

# Organisation website changed from stat.kg to stat.gov.kg.
$ws.Range("B10").Value = "www.stat.gov.kg"

# Both edited cells get their own (new) font/style entry in the workbook.
$ws.Range("B4").Font.Name = "Calibri"
$ws.Range("B10").Font.Name = "Calibri"

# Reflect the author's last selection before saving.
$ws.Range("B6").Select()
